# Updating the model for Horeco
# Shift every timestamp in column A (rows 2..97) forward by one day, and
# update the "Actual Production (MW)" values in column B for rows 2..30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift timestamps in A2:A97 forward by one day ---------------------
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 1
}

# --- New Actual Production (MW) values for rows 2..30 -------------------
$bValues = @(2475,2473,2443,2415,2385,2373,2382,2392,0,2386,2381,2370,2350,2339,2296,2259,2238,2185,2094,1956,1877,1851,1847,1812,1676,1692,1651,1611,1538)

$row = 2
foreach ($val in $bValues) {
    $ws.Cells.Item($row, 2).Value2 = $val
    $row++
}
